$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152; existing rows 152..227 shift down to 153..228.
$ws.Rows(152).Insert()

# Populate the newly inserted row 152 with the new weekly record.
$ws.Cells.Item(152, 1).Value2  = 10
$ws.Cells.Item(152, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value2  = "La Araucanía"
$ws.Cells.Item(152, 4).Value2  = 44452
$ws.Cells.Item(152, 5).Value2  = 9
$ws.Cells.Item(152, 6).Value2  = 100112032
$ws.Cells.Item(152, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(152, 8).Value2  = "Sin especificar"
$ws.Cells.Item(152, 9).Value2  = "Primera"
$ws.Cells.Item(152, 10).Value2 = 500
$ws.Cells.Item(152, 11).Value2 = 17000
$ws.Cells.Item(152, 12).Value2 = 18000
$ws.Cells.Item(152, 13).Value2 = 17600
$ws.Cells.Item(152, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(152, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(152, 16).Value2 = 293
$ws.Cells.Item(152, 17).Value2 = 60
$ws.Cells.Item(152, 18).Value2 = "Hortaliza"
